# Modèle-audit-SEO: add a title banner row above the existing table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new row at the top; everything else shifts down by one. ---
$ws.Rows(1).Insert()

# --- 2. New title row formatting/content. ---
$ws.Range("A1:F1").RowHeight = 40
$ws.Range("A1").Font.Size = 22
$ws.Range("A1:F1").HorizontalAlignment = -4108   # xlCenter
$ws.Range("A1:F1").VerticalAlignment = -4108     # xlCenter
$ws.Range("A1").Value = "le Rapport d'analyse SEO de La Panthère "
$ws.Range("A1:F1").Merge()

# --- 3. Column widths (narrower layout to fit the new banner). ---
$mdwAdjust = 5/7
$ws.Columns("A").ColumnWidth = 13.28515625 - $mdwAdjust
$ws.Columns("B").ColumnWidth = 14.28515625 - $mdwAdjust
$ws.Columns("C").ColumnWidth = 24 - $mdwAdjust
$ws.Columns("D").ColumnWidth = 17.5703125 - $mdwAdjust
$ws.Columns("E").ColumnWidth = 0.140625 - $mdwAdjust
$ws.Columns("F").ColumnWidth = 24.28515625 - $mdwAdjust

# --- 4. Row heights for the header + data rows (now shifted to rows 2-17). ---
$ws.Rows(2).RowHeight = 35
$ws.Rows(3).RowHeight = 75
$ws.Rows(4).RowHeight = 61
$ws.Rows(5).RowHeight = 63
$ws.Rows(6).RowHeight = 66
$ws.Rows(7).RowHeight = 93
$ws.Rows(8).RowHeight = 69
$ws.Rows(9).RowHeight = 64
$ws.Rows(10).RowHeight = 67
$ws.Rows(11).RowHeight = 54
$ws.Rows(12).RowHeight = 64
$ws.Rows(13).RowHeight = 51
$ws.Rows(14).RowHeight = 61
$ws.Rows(15).RowHeight = 78
$ws.Rows(16).RowHeight = 67
$ws.Rows(17).RowHeight = 38

# --- 5. Hyperlinks: the insert doesn't auto-shift hyperlink ranges, so redo them one row down. ---
$links = New-Object System.Collections.ArrayList
foreach ($hl in $ws.Hyperlinks) {
    [void]$links.Add(@($hl.Range.Row, $hl.Range.Column, $hl.Address, $hl.SubAddress))
}
$ws.Hyperlinks.Delete()
foreach ($l in $links) {
    $target = $ws.Cells.Item([int]$l[0] + 1, [int]$l[1])
    $ws.Hyperlinks.Add($target, $l[2], $l[3])
}

# --- 6. Selection state matches the new banner range. ---
$ws.Range("A1:F1").Select()
